$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The daily league-base refresh (29-03-2024) drops the stale/duplicate
# fixture that used to sit on row 136 (id 134, match id 6788938) and
# shifts the remaining upcoming fixtures up by one row. Their odds
# (columns K:V) are also refreshed with the latest values pulled for
# the update.

$ws.Rows(136).Delete()

# --- Fix up the sequential "id" column (A) for the rows that shifted up ---
$ws.Range("A136").Value = 134
$ws.Range("A137").Value = 135
$ws.Range("A138").Value = 136

# --- Row 136 (match id 6788936): refresh odds columns N, O, P, R, S ---
$ws.Range("N136").Value = 2.2
$ws.Range("O136").Value = 3.25
$ws.Range("P136").Value = 3.1
$ws.Range("R136").Value = 1.925
$ws.Range("S136").Value = 1.925

# --- Row 137 (match id 6769306): refresh odds columns R, S, U ---
$ws.Range("R137").Value = 2
$ws.Range("S137").Value = 1.85
$ws.Range("U137").Value = 1.8
